$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update case count for year 2002 from 20 to 19
$ws.Range("B2").Value = 19

# Update the active selection to F18
$ws.Range("F18").Select()
